# Alpha.xlsx geocoding sheet update:
#  - append 3 new geocoded rows to the bottom of the table
#  - re-sort the data range (A2:E.. ) by Location (column A) ascending
#  - left-align columns B:D (keeping the existing wrap-text look on the one
#    row that already had it)
#  - point the saved-file absolute path at the new project folder

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. append the new rows under the current last row (row 15) ----------
$ws.Range("A16").Value = "264 North 10th Street"
$ws.Range("E16").Value = "Can not find the address, not specific"

$ws.Range("A17").Value = "81-111 Junius Street"
$ws.Range("B17").Value = 40.672759999999997
$ws.Range("C17").Value = -73.903649999999999
$ws.Range("D17").Value = "11212-8004"

$ws.Range("A18").Value = "Ralph Avenue and Preston Court"
$ws.Range("B18").Value = 40.641182000000001
$ws.Range("C18").Value = -73.919568999999996
$ws.Range("D18").Value = 11234

# ---- 2. left-align the B:D columns (applies to header + every data row) --
$ws.Range("B1:D1048576").HorizontalAlignment = -4131

# ---- 3. sort A2:E18 by column A, ascending, header excluded --------------
$rng = $ws.Range("A2:E18")
$key = $ws.Range("A2:A18")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($key)
$ws.Sort.SetRange($rng)
$ws.Sort.Header = 0
$ws.Sort.Orientation = 1
$ws.Sort.Apply()

# ---- 4. re-point the stored absolute path ---------------------------------
$wb.Path = "/Users/Zt/Desktop/Geocoding/"
